# Number the "Innovation #:" headings sequentially (Innovation 1, 2, 3, ...)
# in the order they appear in the document.
#
# Each of the six "Significance" sub-headings currently reads literally
# "Innovation #: <title>". We replace just the "#" character with the
# 1-based ordinal of the heading, leaving everything else (text, bold
# run formatting, surrounding paragraph) untouched.

$d = $word.ActiveDocument

$innovationNumber = 0
$paraCount = $d.Paragraphs.Count

for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text

    if ($paraText.StartsWith("Innovation #")) {
        $innovationNumber = $innovationNumber + 1

        $hashOffset = $paraText.IndexOf("#")
        $hashPos = $para.Range.Start + $hashOffset
        $hashRange = $d.Range($hashPos, $hashPos + 1)

        if ($hashRange.Text -eq "#") {
            $hashRange.Text = [string]$innovationNumber
        }
    }
}

Write-Output "Numbered $innovationNumber Innovation heading(s)"
